# Refresh "想去人数" (interested-count) figures on the 展览 and 全部类型
# sheets to match the newly scraped snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1371
$ws1.Range("F7").Value  = 11731
$ws1.Range("F8").Value  = 4397
$ws1.Range("F14").Value = 1100
$ws1.Range("F15").Value = 151
$ws1.Range("F17").Value = 5115
$ws1.Range("F19").Value = 186
$ws1.Range("F20").Value = 518
$ws1.Range("F21").Value = 11351
$ws1.Range("F22").Value = 11289

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1371
$ws4.Range("F7").Value  = 11731
$ws4.Range("F8").Value  = 4397
$ws4.Range("F15").Value = 1100
$ws4.Range("F16").Value = 151
$ws4.Range("F18").Value = 5115
$ws4.Range("F20").Value = 186
$ws4.Range("F21").Value = 518
$ws4.Range("F22").Value = 11351
$ws4.Range("F23").Value = 11289
